$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 2 (RCB)
$ws.Range("B2").Value = 4
$ws.Range("C2").Value = 1
$ws.Range("D2").Value = 3
$ws.Range("E2").Value = 2

# Update row 3 (MI)
$ws.Range("B3").Value = 1
$ws.Range("C3").Value = 1
$ws.Range("D3").Value = 0
$ws.Range("E3").Value = 2

# Update row 5 (DC)
$ws.Range("B5").Value = 3
$ws.Range("C5").Value = 2
$ws.Range("D5").Value = 1
$ws.Range("E5").Value = 4

# Update the selected cell in the sheet view
$ws.Range("G6").Select()
